# Maximum Speed Analysis.docx edit
# Paragraph: "Main assumption in this study is that air enters the chamber at
# Mach 1 which would be the maximum flow condition. Combustion is also assumed
# complete with no cyclic variations."
#
# 1. Insert " For the flow area, intake valve diameter of 34.5 mm from 2026
#    regulations is applied." right after "...maximum flow condition." and
#    before " Combustion is also assumed...".
# 2. Append a trailing space after the final "." of that paragraph.

$d = $word.ActiveDocument

# --- Step 1: insert the new sentence after "maximum flow condition." -------
$rng = $d.Content
$found = $rng.Find.Execute("air enters the chamber at Mach 1 which would be the maximum flow condition.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" For the flow area, ")
    $rng.Collapse(0)
    $rng.InsertAfter("intake")
    $rng.Collapse(0)
    $rng.InsertAfter(" valve diameter of 34.5 mm ")
    $rng.Collapse(0)
    $rng.InsertAfter("from 2026 regulations is applied.")
}

# --- Step 2: append a trailing space after the paragraph's final period ----
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Combustion is also assumed complete with no cyclic variations.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.InsertAfter(" ")
}
